# Swap the values of columns A, B, D, E, F, G, H, Q, R between row 11 and row 12.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell11 = $ws.Range($col + "11")
    $cell12 = $ws.Range($col + "12")

    $v11 = $cell11.Value2
    $v12 = $cell12.Value2

    $cell11.Value2 = $v12
    $cell12.Value2 = $v11
}
